$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 458045
$ws.Range("D2").Value = 536251
$ws.Range("E2").Value = 0.539327323050681
$ws.Range("F2").Value = 185156
$ws.Range("G2").Value = 290733
$ws.Range("H2").Value = 323837
$ws.Range("I2").Value = 462396
$ws.Range("J2").Value = 409734
$ws.Range("K2").Value = 608924
$ws.Range("L2").Value = 558265
$ws.Range("M2").Value = 511379
$ws.Range("N2").Value = 446867
$ws.Range("O2").Value = 468150
$ws.Range("P2").Value = 515830
$ws.Range("Q2").Value = 430808

$ws.Range("C3").Value = 591451
$ws.Range("D3").Value = 583909
$ws.Range("E3").Value = 0.496791621290498
$ws.Range("F3").Value = 242509
$ws.Range("G3").Value = 251120
$ws.Range("H3").Value = 390294
$ws.Range("I3").Value = 613627
$ws.Range("J3").Value = 572019
$ws.Range("K3").Value = 661001
$ws.Range("L3").Value = 651920
$ws.Range("M3").Value = 644582
$ws.Range("N3").Value = 437550
$ws.Range("O3").Value = 459514
$ws.Range("P3").Value = 423504
$ws.Range("Q3").Value = 417638

$ws.Range("C4").Value = 259157
$ws.Range("D4").Value = 147792
$ws.Range("E4").Value = 0.363170815016132
$ws.Range("F4").Value = 63358
$ws.Range("G4").Value = 69640
$ws.Range("H4").Value = 94222
$ws.Range("I4").Value = 113145
$ws.Range("J4").Value = 126842
$ws.Range("K4").Value = 166107
$ws.Range("L4").Value = 179046
$ws.Range("M4").Value = 168835
$ws.Range("N4").Value = 142410
$ws.Range("O4").Value = 162671
$ws.Range("P4").Value = 174984
$ws.Range("Q4").Value = 117911

$ws.Range("C5").Value = 86518
$ws.Range("D5").Value = 155199
$ws.Range("E5").Value = 0.642069031139721
$ws.Range("F5").Value = 67547
$ws.Range("G5").Value = 72975
$ws.Range("H5").Value = 76812
$ws.Range("I5").Value = 96812
$ws.Range("J5").Value = 99266
$ws.Range("K5").Value = 101066
$ws.Range("L5").Value = 110059
$ws.Range("M5").Value = 109981
$ws.Range("N5").Value = 80538
$ws.Range("O5").Value = 105910
$ws.Range("P5").Value = 115116
$ws.Range("Q5").Value = 110134

$ws.Range("C6").Value = 46257
$ws.Range("D6").Value = 36443
$ws.Range("E6").Value = 0.440665054413543
$ws.Range("F6").Value = 20600
$ws.Range("G6").Value = 19813
$ws.Range("H6").Value = 23293
$ws.Range("I6").Value = 29303
$ws.Range("J6").Value = 24225
$ws.Range("K6").Value = 28839
$ws.Range("L6").Value = 24888
$ws.Range("M6").Value = 29851
$ws.Range("N6").Value = 28230
$ws.Range("O6").Value = 41474
$ws.Range("P6").Value = 50455
$ws.Range("Q6").Value = 37957

$ws.Range("C7").Value = 3861486
$ws.Range("D7").Value = 3951895
$ws.Range("E7").Value = 0.505785523578077
$ws.Range("F7").Value = 1486840
$ws.Range("G7").Value = 1848104
$ws.Range("H7").Value = 2340809
$ws.Range("I7").Value = 3600733
$ws.Range("J7").Value = 3459848
$ws.Range("K7").Value = 4396179
$ws.Range("L7").Value = 4348399
$ws.Range("M7").Value = 3901160
$ws.Range("N7").Value = 3287373
$ws.Range("O7").Value = 3405689
$ws.Range("P7").Value = 3439490
$ws.Range("Q7").Value = 2855092

$ws.Range("C8").Value = 85379
$ws.Range("D8").Value = 75496
$ws.Range("E8").Value = 0.469283605283605
$ws.Range("F8").Value = 35190
$ws.Range("G8").Value = 34540
$ws.Range("H8").Value = 59825
$ws.Range("I8").Value = 77186
$ws.Range("J8").Value = 78943
$ws.Range("K8").Value = 87904
$ws.Range("L8").Value = 89111
$ws.Range("M8").Value = 88090
$ws.Range("N8").Value = 58170
$ws.Range("O8").Value = 80316
$ws.Range("P8").Value = 66370
$ws.Range("Q8").Value = 43849

$ws.Range("C9").Value = 273739
$ws.Range("D9").Value = 261011
$ws.Range("E9").Value = 0.488099111734455
$ws.Range("F9").Value = 64255
$ws.Range("G9").Value = 126063
$ws.Range("H9").Value = 139654
$ws.Range("I9").Value = 256531
$ws.Range("J9").Value = 210765
$ws.Range("K9").Value = 342759
$ws.Range("L9").Value = 352324
$ws.Range("M9").Value = 304389
$ws.Range("N9").Value = 265583
$ws.Range("O9").Value = 275799
$ws.Range("P9").Value = 283209
$ws.Range("Q9").Value = 217449

$ws.Range("C10").Value = 64196
$ws.Range("D10").Value = 74800
$ws.Range("E10").Value = 0.538144982589427
$ws.Range("F10").Value = 33035
$ws.Range("G10").Value = 35352
$ws.Range("H10").Value = 50236
$ws.Range("I10").Value = 73825
$ws.Range("J10").Value = 86573
$ws.Range("K10").Value = 97018
$ws.Range("L10").Value = 90804
$ws.Range("M10").Value = 67932
$ws.Range("N10").Value = 61667
$ws.Range("O10").Value = 64944
$ws.Range("P10").Value = 56491
$ws.Range("Q10").Value = 42184

$ws.Range("C11").Value = 357614
$ws.Range("D11").Value = 291941
$ws.Range("E11").Value = 0.449447698809185
$ws.Range("F11").Value = 105526
$ws.Range("G11").Value = 142331
$ws.Range("H11").Value = 229349
$ws.Range("I11").Value = 320210
$ws.Range("J11").Value = 349311
$ws.Range("K11").Value = 452224
$ws.Range("L11").Value = 435544
$ws.Range("M11").Value = 348782
$ws.Range("N11").Value = 229134
$ws.Range("O11").Value = 292509
$ws.Range("P11").Value = 390316
$ws.Range("Q11").Value = 323980

$ws.Range("C12").Value = 870980
$ws.Range("D12").Value = 665944
$ws.Range("E12").Value = 0.433296636658677
$ws.Range("F12").Value = 219054
$ws.Range("G12").Value = 372374
$ws.Range("H12").Value = 479181
$ws.Range("I12").Value = 842494
$ws.Range("J12").Value = 858013
$ws.Range("K12").Value = 1104375
$ws.Range("L12").Value = 1064548
$ws.Range("M12").Value = 900649
$ws.Range("N12").Value = 866340
$ws.Range("O12").Value = 737700
$ws.Range("P12").Value = 590426
$ws.Range("Q12").Value = 506521

$ws.Range("C13").Value = 190656
$ws.Range("D13").Value = 260684
$ws.Range("E13").Value = 0.577577879204148
$ws.Range("F13").Value = 119061
$ws.Range("G13").Value = 96894
$ws.Range("H13").Value = 119735
$ws.Range("I13").Value = 169228
$ws.Range("J13").Value = 140830
$ws.Range("K13").Value = 183654
$ws.Range("L13").Value = 185013
$ws.Range("M13").Value = 189415
$ws.Range("N13").Value = 165371
$ws.Range("O13").Value = 201697
$ws.Range("P13").Value = 215206
$ws.Range("Q13").Value = 153798

$ws.Range("C14").Value = 49812
$ws.Range("D14").Value = 45722
$ws.Range("E14").Value = 0.478594008415852
$ws.Range("F14").Value = 12156
$ws.Range("G14").Value = 22550
$ws.Range("H14").Value = 41550
$ws.Range("I14").Value = 54428
$ws.Range("J14").Value = 48911
$ws.Range("K14").Value = 62467
$ws.Range("L14").Value = 64783
$ws.Range("M14").Value = 48808
$ws.Range("N14").Value = 27508
$ws.Range("O14").Value = 29600
$ws.Range("P14").Value = 38360
$ws.Range("Q14").Value = 32292

$ws.Range("C15").Value = 140214
$ws.Range("D15").Value = 383221
$ws.Range("E15").Value = 0.732127198219454
$ws.Range("F15").Value = 135866
$ws.Range("G15").Value = 118645
$ws.Range("H15").Value = 128733
$ws.Range("I15").Value = 177300
$ws.Range("J15").Value = 159942
$ws.Range("K15").Value = 154796
$ws.Range("L15").Value = 190584
$ws.Range("M15").Value = 169080
$ws.Range("N15").Value = 150651
$ws.Range("O15").Value = 132999
$ws.Range("P15").Value = 161387
$ws.Range("Q15").Value = 140107

$ws.Range("C16").Value = 85701
$ws.Range("D16").Value = 116034
$ws.Range("E16").Value = 0.575180310803777
$ws.Range("F16").Value = 45797
$ws.Range("G16").Value = 56371
$ws.Range("H16").Value = 55026
$ws.Range("I16").Value = 92122
$ws.Range("J16").Value = 93610
$ws.Range("K16").Value = 113064
$ws.Range("L16").Value = 108315
$ws.Range("M16").Value = 90497
$ws.Range("N16").Value = 84366
$ws.Range("O16").Value = 86517
$ws.Range("P16").Value = 94540
$ws.Range("Q16").Value = 82424

$ws.Range("C17").Value = 177262
$ws.Range("D17").Value = 143903
$ws.Range("E17").Value = 0.448065636043778
$ws.Range("F17").Value = 75393
$ws.Range("G17").Value = 85635
$ws.Range("H17").Value = 66139
$ws.Range("I17").Value = 137359
$ws.Range("J17").Value = 108772
$ws.Range("K17").Value = 132756
$ws.Range("L17").Value = 143716
$ws.Range("M17").Value = 128558
$ws.Range("N17").Value = 130025
$ws.Range("O17").Value = 142380
$ws.Range("P17").Value = 136961
$ws.Range("Q17").Value = 100283

$ws.Range("C18").Value = 122906
$ws.Range("D18").Value = 166993
$ws.Range("E18").Value = 0.576038551357542
$ws.Range("F18").Value = 56492
$ws.Range("G18").Value = 46521
$ws.Range("H18").Value = 49645
$ws.Range("I18").Value = 71887
$ws.Range("J18").Value = 78216
$ws.Range("K18").Value = 86978
$ws.Range("L18").Value = 87912
$ws.Range("M18").Value = 86605
$ws.Range("N18").Value = 98585
$ws.Range("O18").Value = 99046
$ws.Range("P18").Value = 108498
$ws.Range("Q18").Value = 84266
